$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "57.730.24"
Set-TextValue $ws.Range("E2") "  +0.12%  "
Set-TextValue $ws.Range("D3") "3.062.49"
Set-TextValue $ws.Range("E3") "  +1.71%  "
Set-TextValue $ws.Range("E4") "  +0.03%  "
Set-TextValue $ws.Range("D5") "517.43"
Set-TextValue $ws.Range("E5") "  +1.32%  "
Set-TextValue $ws.Range("D6") "140.41"
Set-TextValue $ws.Range("E6") "  +0.63%  "
Set-TextValue $ws.Range("E7") "  +0.06%  "
Set-TextValue $ws.Range("D8") "0.434"
Set-TextValue $ws.Range("E8") "  -0.73%  "
Set-TextValue $ws.Range("D9") "7.30"
Set-TextValue $ws.Range("E9") "  -3.69%  "
Set-TextValue $ws.Range("E10") "  -0.57%  "
Set-TextValue $ws.Range("D11") "0.376"
Set-TextValue $ws.Range("E11") "  +2.97%  "
Set-TextValue $ws.Range("D12") "3.578.26"
Set-TextValue $ws.Range("E12") "  +1.49%  "
Set-TextValue $ws.Range("E13") "  -3.22%  "
Set-TextValue $ws.Range("D14") "26.85"
Set-TextValue $ws.Range("E14") "  +1.38%  "
Set-TextValue $ws.Range("D15") "0.0000168"
Set-TextValue $ws.Range("E15") "  +2.74%  "
Set-TextValue $ws.Range("D16") "57.568.94"
Set-TextValue $ws.Range("E16") "  -0.09%  "
Set-TextValue $ws.Range("D17") "6.23"
Set-TextValue $ws.Range("E17") "  +0.44%  "
Set-TextValue $ws.Range("D18") "3.066.32"
Set-TextValue $ws.Range("E18") "  +1.80%  "
Set-TextValue $ws.Range("D19") "13.39"
Set-TextValue $ws.Range("E19") "  +4.23%  "
Set-TextValue $ws.Range("D20") "8.21"
Set-TextValue $ws.Range("E20") "  +3.00%  "
Set-TextValue $ws.Range("D21") "330.61"
Set-TextValue $ws.Range("E21") "  -0.17%  "
Set-TextValue $ws.Range("E22") "  +0.14%  "
Set-TextValue $ws.Range("D23") "0.508"
Set-TextValue $ws.Range("E23") "  +2.03%  "
Set-TextValue $ws.Range("D24") "66.03"
Set-TextValue $ws.Range("E24") "  +2.22%  "
Set-TextValue $ws.Range("D25") "3.175.54"
Set-TextValue $ws.Range("E25") "  +1.25%  "
Set-TextValue $ws.Range("E26") "  -2.55%  "
Set-TextValue $ws.Range("E27") "  -0.11%  "
Set-TextValue $ws.Range("D28") "0.0₃0907"
Set-TextValue $ws.Range("E28") "  -1.74%  "
Set-TextValue $ws.Range("D29") "6.73"
Set-TextValue $ws.Range("E29") "  -0.90%  "
Set-TextValue $ws.Range("D30") "7.32"
Set-TextValue $ws.Range("E30") "  -0.51%  "
Set-TextValue $ws.Range("E31") "  -0.40%  "
Set-TextValue $ws.Range("E32") "  +1.85%  "
Set-TextValue $ws.Range("D33") "20.87"
Set-TextValue $ws.Range("E33") "  +1.45%  "
Set-TextValue $ws.Range("D34") "153.42"
Set-TextValue $ws.Range("E34") "  -0.90%  "
Set-TextValue $ws.Range("D35") "4.65"
Set-TextValue $ws.Range("E35") "  -1.93%  "
Set-TextValue $ws.Range("D36") "5.91"
Set-TextValue $ws.Range("E36") "  +0.59%  "
Set-TextValue $ws.Range("D37") "25.55"
Set-TextValue $ws.Range("E37") "  +4.65%  "
Set-TextValue $ws.Range("D38") "1.28"
Set-TextValue $ws.Range("E38") "  +0.05%  "
Set-TextValue $ws.Range("D39") "0.0679"
Set-TextValue $ws.Range("E39") "  +0.38%  "
Set-TextValue $ws.Range("D40") "37.09"
Set-TextValue $ws.Range("E40") "  -1.14%  "
Set-TextValue $ws.Range("D41") "3.89"
Set-TextValue $ws.Range("E41") "  +0.86%  "
Set-TextValue $ws.Range("D42") "0.670"
Set-TextValue $ws.Range("E42") "  +2.96%  "
Set-TextValue $ws.Range("D43") "0.998"
Set-TextValue $ws.Range("E43") "  -0.22%  "
Set-TextValue $ws.Range("B44") "Stacks"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D44") "1.40"
Set-TextValue $ws.Range("E44") "  -1.00%  "
Set-TextValue $ws.Range("B45") "Maker"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D45") "2.205.53"
Set-TextValue $ws.Range("E45") "  -1.15%  "
Set-TextValue $ws.Range("D46") "6.13"
Set-TextValue $ws.Range("E46") "  +1.89%  "
Set-TextValue $ws.Range("E47") "  +2.33%  "
Set-TextValue $ws.Range("D48") "0.956"
Set-TextValue $ws.Range("E48") "  -2.90%  "
Set-TextValue $ws.Range("D49") "20.12"
Set-TextValue $ws.Range("E49") "  +3.58%  "
Set-TextValue $ws.Range("E50") "  -5.09%  "
Set-TextValue $ws.Range("B51") "Notcoin"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/2L2Y4ghjj+notcoin-not"
Set-TextValue $ws.Range("D51") "0.0172"
Set-TextValue $ws.Range("E51") "  +10.04%  "
